{"js": "// Word JS API (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Target edit (first paragraph of the document, the \"**ID__..._ID**\" tag\n// paragraph):\n//   1. Add a paragraph border (top/left/bottom/right) whose only attribute\n//      is a 5-twip \"space\" (distance from text) \u2014 no line/color.\n//   2. Change the paragraph's left indent from 120 -> 225 twips\n//      (= 6pt -> 11.25pt).\n//   3. Replace the run text \"**ID__AFFARS_5343_topic_5__ID**\" with\n//      \"**ID__AFFARS_SUBPART_5343_2__ID**\" and drop the trailing\n//      \" \" (space) run that followed it, leaving a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// --- 1 & 3: rewrite the paragraph's text as a single run -----------------\n// Replacing the whole-paragraph range's text collapses every run in the\n// paragraph into one run that carries the first run's formatting (rFonts /\n// b / i / color / sz), which both updates the id text AND removes the\n// trailing \" \" run in one step.\nconst paragraphRange = firstParagraph.getRange();\nparagraphRange.insertText(\"**ID__AFFARS_SUBPART_5343_2__ID**\", \"Replace\");\nawait context.sync();\n\n// --- 2: left indent 120 -> 225 twips (11.25pt) ----------------------------\nfirstParagraph.paragraphFormat.leftIndent = 11.25;\nawait context.sync();\n\n// --- 3: paragraph border (<w:pBdr><w:top w:space=\"5\"/>... ) --------------\n// Word's JS API `ParagraphBorder` (paragraph.borders items) only exposes\n// type/color/width \u2014 it has no \"space\"/distanceFromText setter \u2014 so the\n// border spacing is written through the same low-level OM bridge the\n// Word.js shim itself is built on (`__native.docxOmSet`, mirroring\n// `Borders.DistanceFromTop/Left/Bottom/Right` from the Word object model).\nconst h = firstParagraph._h;\nconst a = firstParagraph._a;\nfor (const side of [\"Top\", \"Left\", \"Bottom\", \"Right\"]) {\n  __native.docxOmSet(h, a, `Borders.DistanceFrom${side}`, \"5\");\n}\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Target edit (first paragraph of the document, the \"**ID__..._ID**\" tag\n# paragraph):\n#   1. Add a paragraph border (top/left/bottom/right) whose only attribute\n#      is a 5-twip \"space\" (distance from text) \u2014 no line/color.\n#   2. Change the paragraph's left indent from 120 -> 225 twips\n#      (= 6pt -> 11.25pt).\n#   3. Replace the run text \"**ID__AFFARS_5343_topic_5__ID**\" with\n#      \"**ID__AFFARS_SUBPART_5343_2__ID**\" and drop the trailing\n#      \" \" (space) run that followed it, leaving a single run.\n\n$d = $word.ActiveDocument\n$p1 = $d.Paragraphs.Item(1)\n\n# --- 1 & 3: rewrite the paragraph's text as a single run -----------------\n# Assigning the whole-paragraph Range.Text collapses every run in the\n# paragraph into one run that carries the first run's formatting (rFonts /\n# b / i / color / sz), which both updates the id text AND removes the\n# trailing \" \" run in one step.\n$p1.Range.Text = \"**ID__AFFARS_SUBPART_5343_2__ID**\"\n\n# --- 2: left indent 120 -> 225 twips (11.25pt) ----------------------------\n$p1.LeftIndent = 11.25\n\n# --- 3: paragraph border (<w:pBdr><w:top w:space=\"5\"/>... ) --------------\n$p1.Borders.DistanceFromTop = 5\n$p1.Borders.DistanceFromLeft = 5\n$p1.Borders.DistanceFromBottom = 5\n$p1.Borders.DistanceFromRight = 5\n"}
